# Update the "HP values" material lookup table in the "Lookup Tables" sheet,
# and repoint the Calculations sheet at the new, larger table.

$wb = $excel.ActiveWorkbook
$wsCalc = $wb.Worksheets.Item("Calculations")
$wsLookup = $wb.Worksheets.Item("Lookup Tables")

# --- Lookup Tables!D1:E18 : material name / HP factor table -----------------
# Rows 1-9 : existing materials, re-cased to title case, HP values refreshed
$wsLookup.Range("D1").Value = "Aluminum"
$wsLookup.Range("E1").Value = 0.25

$wsLookup.Range("D2").Value = "Magnesium"
$wsLookup.Range("E2").Value = 0.25

$wsLookup.Range("D3").Value = "Copper"
$wsLookup.Range("E3").Value = 0.5

$wsLookup.Range("D4").Value = "Brass"
$wsLookup.Range("E4").Value = 0.4

$wsLookup.Range("D5").Value = "Bronze"
$wsLookup.Range("E5").Value = 0.5

$wsLookup.Range("D6").Value = "Steel(up to 150 Brinell)"
$wsLookup.Range("E6").Value = 1.4

$wsLookup.Range("D7").Value = "Steel(up to 300 Brinell)"
$wsLookup.Range("E7").Value = 1.7

$wsLookup.Range("D8").Value = "Steel(up to 400 Brinell)"
$wsLookup.Range("E8").Value = 2

$wsLookup.Range("D9").Value = "Steel(up to 500 Brinell)"
$wsLookup.Range("E9").Value = 2.5

# Rows 10-18 : brand new materials appended to the table
$wsLookup.Range("D10").Value = "Gray Cast Iron"
$wsLookup.Range("E10").Value = 0.5

$wsLookup.Range("D11").Value = "Ductile Cast Iron"
$wsLookup.Range("E11").Value = 0.56

$wsLookup.Range("D12").Value = "Maleable Cast Iron"
$wsLookup.Range("E12").Value = 0.67

$wsLookup.Range("D13").Value = "Chilled Cast Iron"
$wsLookup.Range("E13").Value = 1.67

$wsLookup.Range("D14").Value = "High Tensile Alloys"
$wsLookup.Range("E14").Value = 2.5

$wsLookup.Range("D15").Value = "Titanium"
$wsLookup.Range("E15").Value = 1.67

$wsLookup.Range("D16").Value = "PH Series Stainless Steels"
$wsLookup.Range("E16").Value = 1.33

$wsLookup.Range("D17").Value = "300 Series Stainless Steels"
$wsLookup.Range("E17").Value = 1

$wsLookup.Range("D18").Value = "High Temp Alloys"
$wsLookup.Range("E18").Value = 2.5

# Move the VLOOKUP helper formula down to row 19 and widen its table range
$wsLookup.Range("E19").Formula = "=VLOOKUP(Calculations!B4,D1:E18,2,0)"

# Widen column D a bit to fit the longer material names
$wsLookup.Columns("D").ColumnWidth = 24.85546875

# --- Calculations sheet ------------------------------------------------------
# Re-select the Aluminum material (now title-cased) in the (now larger) list
$wsCalc.Range("B4").Value = "Aluminum"

# Point the Cutting Power (HP) formula at the relocated lookup cell
$wsCalc.Range("B23").Formula = "=B22*B9*B20*'Lookup Tables'!E19"

# --- Data validation on Calculations!B4 -------------------------------------
# Extend the dropdown source range to cover the new rows, and stop disabling
# the input prompts (Excel re-enables them whenever the validation list is
# edited through the UI).
$dv = $wsCalc.Range("B4").Validation
$dv.Delete()
$dv.Add(3, 1, 1, "='Lookup Tables'!`$D`$1:`$D`$18")
$dv.IgnoreBlank = $true
$dv.InCellDropdown = $true
$dv.ShowInput = $true
$dv.ShowError = $true
